$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-51 with refreshed crypto data.
# Values that look like plain numbers would otherwise be auto-converted from text
# to a Number by Excel, so those cells are forced to Text (matching the source data,
# which stores these as inline strings) and then restored to the default "Normal"
# cell style so no stray formatting is introduced.

$ws.Range("D2").Value = "27.392.89"
$ws.Range("E2").Value = "  +4.72%  "

$ws.Range("D3").Value = "1.817.14"
$ws.Range("E3").Value = "  +5.67%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.70%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "344.78"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.11%  "

$ws.Range("E6").Value = "  +0.43%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3814"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.14%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3512"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +4.55%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "49.44"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.28%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.231"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.60%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07751"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.50%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.80%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "22.12"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +10.11%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.627"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +5.55%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.273"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +4.91%  "

$ws.Range("D16").Value = "1.817.35"
$ws.Range("E16").Value = "  +5.62%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001123"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +4.02%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.06733"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.04%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "86.23"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.95%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.40%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.69"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +7.75%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.541"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +7.37%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "13.23"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.62%  "

$ws.Range("D24").Value = "27.414.10"
$ws.Range("E24").Value = "  +5.17%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.468"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.675"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +7.17%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "22.14"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +14.50%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.486"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +12.31%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "154.30"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.76%  "

$ws.Range("D30").Value = "2.025.08"
$ws.Range("E30").Value = "  +6.16%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "136.16"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +5.22%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.333"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.75%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.052"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.21%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "13.93"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +7.34%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.08763"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.85%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.706"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.74%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.639"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.80%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.7025"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +13.48%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.2278"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +6.58%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.06529"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.82%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.02411"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +5.13%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "9.018"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +5.40%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.303"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.44%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "14.75"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.28%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.6553"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +11.06%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.39%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.028"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.93%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.184"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +8.24%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "132.63"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.80%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.07343"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.86%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "80.76"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.58%  "
